$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change F2 from text "20201005" to the number 20201005
$ws.Range("F2").Value = 20201005

# Move the active selection to H12
$ws.Range("H12").Select()
